$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 418.54544
$ws.Range("I2").Value = 452
$ws.Range("K2").Value = 452
$ws.Range("M2").Value = -339
$ws.Range("H17").Value = 1421758
$ws.Range("J17").Value = 1516455.2
$ws.Range("L17").Value = 4549365.6
$ws.Range("N17").Value = -4549701.6
$ws.Range("H100").Value = 165900.44
$ws.Range("I100").Value = 200660.8
$ws.Range("K100").Value = 200660.8
$ws.Range("M100").Value = -200119.8
$ws.Range("H112").Value = 2587.2104
$ws.Range("J112").Value = 2223.9412
$ws.Range("L112").Value = 6671.823600000001
$ws.Range("N112").Value = -8887.8236
$ws.Range("H132").Value = 2861383.8
$ws.Range("J132").Value = 33338666
$ws.Range("L132").Value = 100015998
$ws.Range("N132").Value = -100021058
$ws.Range("H133").Value = 91098.60000000001
$ws.Range("I133").Value = 70709
$ws.Range("J133").Value = 93364.11
$ws.Range("K133").Value = 70709
$ws.Range("L133").Value = 93364.11
$ws.Range("M133").Value = -65649
$ws.Range("N133").Value = -103484.11
$ws.Range("H136").Value = 67999.336
$ws.Range("J136").Value = 67999.336
$ws.Range("L136").Value = 67999.336
$ws.Range("N136").Value = -78199.336
$ws.Range("H138").Value = 357868.62
$ws.Range("I138").Value = 549540.0600000001
$ws.Range("K138").Value = 1648620.18
$ws.Range("M138").Value = -1643480.18
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 10701.2
$ws.Range("I61").Value = 14503.308
$ws.Range("K61").Value = 14503.308
$ws.Range("M61").Value = -14291.308
$ws.Range("H132").Value = 3846.2917
$ws.Range("I132").Value = 2108.2307
$ws.Range("K132").Value = 6324.6921
$ws.Range("M132").Value = -3794.6921
$ws.Range("H135").Value = 97987.734
$ws.Range("J135").Value = 97987.734
$ws.Range("L135").Value = 97987.734
$ws.Range("N135").Value = -108127.734
$ws.Range("H136").Value = 10701.2
$ws.Range("I136").Value = 14503.308
$ws.Range("K136").Value = 43509.924
$ws.Range("M136").Value = -40959.924
$ws.Range("H138").Value = 62027.5
$ws.Range("J138").Value = 62027.5
$ws.Range("L138").Value = 62027.5
$ws.Range("N138").Value = -72307.5
$ws.Range("H141").Value = 49835.832
$ws.Range("J141").Value = 49835.832
$ws.Range("L141").Value = 49835.832
$ws.Range("N141").Value = -60195.832
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1019666.3
$ws.Range("I4").Value = 9999
$ws.Range("J4").Value = 1524500
$ws.Range("K4").Value = 9999
$ws.Range("L4").Value = 1524500
$ws.Range("M4").Value = -9887
$ws.Range("N4").Value = -1524724
$ws.Range("H16").Value = 622.7826
$ws.Range("I16").Value = 670.1053000000001
$ws.Range("K16").Value = 670.1053000000001
$ws.Range("M16").Value = -383.1053000000001
$ws.Range("H31").Value = 12538.615
$ws.Range("I31").Value = 19286.143
$ws.Range("J31").Value = 4666.5
$ws.Range("K31").Value = 19286.143
$ws.Range("L31").Value = 4666.5
$ws.Range("M31").Value = -18991.143
$ws.Range("N31").Value = -5256.5
$ws.Range("H34").Value = 12538.615
$ws.Range("I34").Value = 19286.143
$ws.Range("J34").Value = 4666.5
$ws.Range("K34").Value = 19286.143
$ws.Range("L34").Value = 4666.5
$ws.Range("M34").Value = -19084.143
$ws.Range("N34").Value = -5070.5
$ws.Range("H86").Value = 6370.8
$ws.Range("I86").Value = 5027.3335
$ws.Range("J86").Value = 8386
$ws.Range("K86").Value = 5027.3335
$ws.Range("L86").Value = 8386
$ws.Range("M86").Value = -3904.3335
$ws.Range("N86").Value = -10632
$ws.Range("H89").Value = 6370.8
$ws.Range("I89").Value = 5027.3335
$ws.Range("J89").Value = 8386
$ws.Range("K89").Value = 25136.6675
$ws.Range("L89").Value = 41930
$ws.Range("M89").Value = -19520.6675
$ws.Range("N89").Value = -53162
$ws.Range("H113").Value = 622.7826
$ws.Range("I113").Value = 670.1053000000001
$ws.Range("K113").Value = 670.1053000000001
$ws.Range("M113").Value = 1499.8947
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6254.5835
$ws.Range("I80").Value = 8300.375
$ws.Range("J80").Value = 2163
$ws.Range("K80").Value = 8300.375
$ws.Range("L80").Value = 2163
$ws.Range("M80").Value = -7302.375
$ws.Range("N80").Value = -4159
$ws.Range("H83").Value = 6254.5835
$ws.Range("I83").Value = 8300.375
$ws.Range("J83").Value = 2163
$ws.Range("K83").Value = 41501.875
$ws.Range("L83").Value = 10815
$ws.Range("M83").Value = -36509.875
$ws.Range("N83").Value = -20799
$ws.Range("H93").Value = 42832.25
$ws.Range("J93").Value = 42832.25
$ws.Range("L93").Value = 42832.25
$ws.Range("N93").Value = -46576.25
$ws.Range("H105").Value = 89333
$ws.Range("J105").Value = 89333
$ws.Range("L105").Value = 89333
$ws.Range("N105").Value = -96321
$ws.Range("H132").Value = 4289.364
$ws.Range("J132").Value = 2685.3333
$ws.Range("L132").Value = 8055.999899999999
$ws.Range("N132").Value = -13115.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4358437.5
$ws.Range("I46").Value = 719.8570999999999
$ws.Range("K46").Value = 719.8570999999999
$ws.Range("M46").Value = -531.8570999999999
$ws.Range("H68").Value = 3380.4783
$ws.Range("I68").Value = 2188.6
$ws.Range("J68").Value = 5615.25
$ws.Range("K68").Value = 2188.6
$ws.Range("L68").Value = 5615.25
$ws.Range("M68").Value = -1439.6
$ws.Range("N68").Value = -7113.25
$ws.Range("H71").Value = 3380.4783
$ws.Range("I71").Value = 2188.6
$ws.Range("J71").Value = 5615.25
$ws.Range("K71").Value = 10943
$ws.Range("L71").Value = 28076.25
$ws.Range("M71").Value = -7199
$ws.Range("N71").Value = -35564.25
$ws.Range("H82").Value = 4324.5713
$ws.Range("I82").Value = 4378.6665
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 4378.6665
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -4017.6665
$ws.Range("N82").Value = -4722
$ws.Range("H85").Value = 4324.5713
$ws.Range("I85").Value = 4378.6665
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 4378.6665
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -3130.6665
$ws.Range("N85").Value = -6496
$ws.Range("H93").Value = 4855.773
$ws.Range("I93").Value = 6654.4287
$ws.Range("J93").Value = 1708.125
$ws.Range("K93").Value = 6654.4287
$ws.Range("L93").Value = 1708.125
$ws.Range("M93").Value = -5406.4287
$ws.Range("N93").Value = -4204.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 3500
$ws.Range("I39").Value = 3500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 3500
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("M39").Value = -3087
$ws.Range("H4").Value = 469.45
$ws.Range("I4").Value = 99.545456
$ws.Range("K4").Value = 99.545456
$ws.Range("M4").Value = 13.454544
$ws.Range("H62").Value = 289918.25
$ws.Range("I62").Value = 569734
$ws.Range("K62").Value = 569734
$ws.Range("M62").Value = -569110
$ws.Range("H65").Value = 289918.25
$ws.Range("I65").Value = 569734
$ws.Range("K65").Value = 2848670
$ws.Range("M65").Value = -2845550
$ws.Range("H100").Value = 30546.154
$ws.Range("I100").Value = 7651
$ws.Range("K100").Value = 15302
$ws.Range("M100").Value = -14761
